$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "A2" = "多氟多";  "B2" = "平潭发展"; "C2" = "合富中国"
    "A3" = "天赐材料"; "B3" = "多氟多";  "C3" = "闻泰科技"
    "A4" = "天际股份"; "B4" = "特变电工"; "C4" = "平潭发展"
    "A5" = "方正电机"; "B5" = "海马汽车"; "C5" = "海马汽车"
    "A6" = "平潭发展"; "B6" = "天赐材料"; "C6" = "多氟多"
    "A7" = "海马汽车"; "B7" = "闻泰科技"; "C7" = "特变电工"
    "A8" = "合富中国"; "B8" = "合富中国"; "C8" = "海陆重工"
    "A9" = "特变电工"; "B9" = "方正电机"; "C9" = "天际股份"
    "A10" = "闻泰科技"; "B10" = "三花智控"; "C10" = "万向钱潮"
    "A11" = "雪人集团"; "B11" = "天际股份"; "C11" = "漳州发展"
    "A12" = "万向钱潮"; "B12" = "万向钱潮"; "C12" = "粤传媒"
    "A13" = "三花智控"; "B13" = "雪人集团"; "C13" = "兰石重装"
    "A14" = "海陆重工"; "B14" = "福龙马";  "C14" = "中毅达"
    "A15" = "兰石重装"; "B15" = "兰石重装"; "C15" = "淳中科技"
    "A16" = "福龙马";  "B16" = "海陆重工"; "C16" = "海天股份"
    "A17" = "中国西电"; "B17" = "中国西电"; "C17" = "海峡创新"
    "A18" = "海新能科"; "B18" = "隆基绿能"; "C18" = "隆基绿能"
    "A19" = "山高环能"; "B19" = "通威股份"; "C19" = "吉视传媒"
    "A20" = "卓越新能"; "B20" = "漳州发展"; "C20" = "神州信息"
    "A21" = "五洲新春"; "B21" = "海南发展"; "C21" = "中国西电"
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
